$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.473.52"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.656.76"
$ws.Range("E3").Value = "  -2.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.97"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3627"
$ws.Range("E7").Value = "  -2.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.29"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3268"
$ws.Range("E9").Value = "  -4.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.124"
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06962"
$ws.Range("E11").Value = "  -6.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9986"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.924"
$ws.Range("E13").Value = "  -4.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.31"
$ws.Range("E14").Value = "  -6.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.614"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("D16").Value = "1.654.37"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001042"
$ws.Range("E17").Value = "  -6.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06517"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9987"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "76.36"
$ws.Range("E20").Value = "  -7.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.908"
$ws.Range("E21").Value = "  -6.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.71"
$ws.Range("E22").Value = "  -7.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.57"
$ws.Range("E23").Value = "  -3.72%  "
$ws.Range("D24").Value = "24.423.00"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.461"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.311"
$ws.Range("E26").Value = "  -16.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.24"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.41"
$ws.Range("E28").Value = "  -8.34%  "
$ws.Range("D29").Value = "1.838.80"
$ws.Range("E29").Value = "  -2.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.193"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.18"
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.061"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.582"
$ws.Range("E33").Value = "  -16.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08361"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("E35").Value = "  -4.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.34"
$ws.Range("E36").Value = "  -8.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.207"
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06048"
$ws.Range("E38").Value = "  -6.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02194"
$ws.Range("E39").Value = "  -7.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.204"
$ws.Range("E40").Value = "  -5.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2048"
$ws.Range("E41").Value = "  -5.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.171"
$ws.Range("E42").Value = "  -7.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9987"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5863"
$ws.Range("E44").Value = "  -8.06%  "
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.66"
$ws.Range("E46").Value = "  -7.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5582"
$ws.Range("E47").Value = "  -7.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.01"
$ws.Range("E48").Value = "  -5.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.940"
$ws.Range("E49").Value = "  -7.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06901"
$ws.Range("E50").Value = "  -4.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.90"
$ws.Range("E51").Value = "  -5.95%  "
